$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TP row (row 2): add +1.8V to the list of test-point nets ---
$ws.Range("B2").Value = "+1.8V, +3.3V, +5V_USB, DONE, FCS, FMISO, FMOSI, FSCK, GND, INTN, JTEN, PGMN, Vref"
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# TP quantity bumped 12 -> 13
$ws.Range("G2").Value = 13

# --- 10V 1uF cap row (row 4): designators now include C31, C32 (RAM decoupling) ---
$ws.Range("B4").Value = "C12, C31, C32"
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# Quantity bumped 1 -> 3
$ws.Range("G4").Value = 3

# --- Ext. Power header (row 16): 3-position -> 4-position connector (BGA clearance fix) ---
$ws.Range("C16").Value = "0.1HDR1X4P"
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("E16").Value = "22-11-2042"
$ws.Range("E15").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("F16").Value = 'Connector, Header, 0.1", 4 Positions, Male, 0.1" Spacing, 4A, 250V, 0.24" Long, 0.125" Tail, 15u" Select Gold, Polarized, Locking Ramp, Vert'
$ws.Range("F15").Copy()
$ws.Range("F16").PasteSpecial(-4122)

$ws.Range("I16").Value = "WM2702-ND"
$ws.Range("I15").Copy()
$ws.Range("I16").PasteSpecial(-4122)

# --- New row 31: 1.8V LDO (MCP1703T-1802E/CB) for RAM supply, U9 ---
$ws.Range("A31").Value = "1.8V"
$ws.Range("B31").Value = "U9"
$ws.Range("C31").Value = "SOT23-12-3TOP"
$ws.Range("D31").Value = "Microchip"
$ws.Range("E31").Value = "MCP1703T-1802E/CB"
$ws.Range("F31").Value = 'IC, Voltage Regulator, 1.8Vout, 2.7-16Vin, 200mA, -40°C ~ 125°C, SOT-23, SMD'
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = "Digikey"
$ws.Range("I31").Value = "MCP1703T-1802E/CBCT-ND"

$ws.Range("A30:I30").Copy()
$ws.Range("A31:I31").PasteSpecial(-4122)
